# Applies the "Requirements Progress" updates described in the commit:
#   "Added support for inserting videos into solutions"
#
# Concretely this:
#  - marks several existing requirement rows as "Done"
#  - marks one requirement row as "In Progress"
#  - rewords two requirement descriptions
#  - appends a brand new completed requirement row ("Add support for adding images")
# on the "Developmnet PB" worksheet (the first/active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Colors used throughout this tracker (standard RGB, encoded as BGR ints for COM)
$Green  = 5287936   # RGB(0,176,80)   -> "Done"
$Red    = 255        # RGB(255,0,0)
$Orange = 49407       # RGB(255,192,0)
$Yellow = 65535       # RGB(255,255,0)

function Set-Done($row) {
    $a = $ws.Cells.Item($row, 1)
    $a.Interior.Color = $Green
    $a.WrapText = $true
    $ws.Cells.Item($row, 3).Value2 = "Done"
}

function Set-Status($row, $status) {
    $ws.Cells.Item($row, 3).Value2 = $status
}

# Row 20: "add login/logout for patient handler or administrator" -> Done
Set-Done 20

# Row 38: "Complete Report" -> Done
Set-Done 38

# Row 42: reword item text (status/priority unchanged)
$ws.Cells.Item(42, 1).Value2 = "Change length of numeric attribute form box to be smaller"

# Row 44: "Get rid of `"Index`" from all of the pages..." -> Done
Set-Done 44

# Row 46: "Find out what browsers are compatible" -> In Progress
Set-Status 46 "In Progress"

# Row 53: "Fix bug where when editing tree..." -> Done
Set-Done 53

# Row 54: give the Priority cell (B54) an orange fill (previously unset)
$ws.Cells.Item(54, 2).Interior.Color = $Orange

# Row 56: reword item text and give Priority cell (B56) an orange fill
$ws.Cells.Item(56, 1).Value2 = "Fix bug with weight, i.e. make that attribute non editable"
$ws.Cells.Item(56, 2).Interior.Color = $Orange

# Row 57 (new): "Add support for adding images" - Done, orange priority
$ws.Cells.Item(57, 1).Value2 = "Add support for adding images"
$ws.Cells.Item(57, 1).Interior.Color = $Green
$ws.Cells.Item(57, 1).WrapText = $true
$ws.Cells.Item(57, 2).Interior.Color = $Orange
$ws.Cells.Item(57, 3).Value2 = "Done"

# Update the view so the active cell / visible area matches the edited region
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A35").Select()

Write-Host "Applied Requirements Progress updates"
